# Battleships armor calc workbook: refactor combat values + cosmetic view tweaks.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Avg. Damage vs. Target")

# --- Row 21: 'Sadist' Particle Beam -----------------------------------
# Base damage (Dmg/shot) bumped from 3 to 8.
$ws1.Range("D21").Value = 8
# Special-ability note switched from the old "locks on target" text to a
# flat shield-pierce penalty; force literal text (leading apostrophe) so
# Excel keeps the quotePrefix + wrap formatting applied to the new string.
$ws1.Range("T21").Value = "'-50% Shield Pierce"

# --- Row 29: 'Beast' Plasma Torpedo ------------------------------------
$ws1.Range("D29").Value = 900
$ws1.Range("H29").Value = 120

# --- Row 36: 'Nova' Antimatter Cannon -----------------------------------
$ws1.Range("D36").Value = 1000
$ws1.Range("E36").Value = 30

# --- Row 51: 'Apocalypse' Giga Laser ------------------------------------
$ws1.Range("D51").Value = 200
$ws1.Range("E51").Value = 1
$ws1.Range("H51").Value = 100

# --- Row 52: 'Terminator' Graviton Devastator ---------------------------
$ws1.Range("D52").Value = 10000
$ws1.Range("E52").Value = 100
$ws1.Range("G52").Value = 30
$ws1.Range("H52").Value = 1000

# --- Row 53: 'Doorkeeper' Tachyon Lance ---------------------------------
$ws1.Range("D53").Value = 230
$ws1.Range("G53").Value = 1000

# --- Cosmetic view state (best effort) ----------------------------------
$ws1.Activate()
$ws1.Range("G29").Select()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 21
$aw.ScrollColumn = 12
